$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 23 - FIFA World Cup
$ws.Range("D23").Value = "https://www.kaggle.com/abecklas/fifa-world-cup"
$ws.Range("A23").Value = "FIFA World Cup"
$ws.Range("B23").Value = "All the results from World cups"
$ws.Range("C23").Value = "WorldCupMatches"

# Rows 24-25 - Student Alcohol Consumption
$ws.Range("C24").Value = "student-por"
$ws.Range("C25").Value = "student-mat"
$ws.Range("D24").Value = "https://www.kaggle.com/uciml/student-alcohol-consumption"
$ws.Range("D25").Value = "https://www.kaggle.com/uciml/student-alcohol-consumption"
$ws.Range("A24").Value = "Student Alcohol Consumption"
$ws.Range("A25").Value = "Student Alcohol Consumption"
$ws.Range("B24").Value = "Social, gender and study data from secondary school students (portugese)"
$ws.Range("B25").Value = "Social, gender and study data from secondary school students (math)"

# Row 26 - Chocolate Bar Ratings
$ws.Range("D26").Value = "https://www.kaggle.com/rtatman/chocolate-bar-ratings"
$ws.Range("C26").Value = "flavors_of_cacao"
$ws.Range("A26").Value = "Chocolate Bar Ratings"
$ws.Range("B26").Value = "Expert ratings of over 1,700 chocolate bars"

$ws.Application.ActiveWindow.ScrollRow = 7
$ws.Range("B26").Select() | Out-Null
